$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114
$ws.Range("A114").Value = 113
$ws.Range("B114").Value = "Monday, Jan 09"
$ws.Range("C114").Value = "1:30 PM"
$ws.Range("D114").Value = "LH1376"
$ws.Range("E114").Value = "Frankfurt"
$ws.Range("F114").Value = "(FRA)"
$ws.Range("G114").Value = "Lufthansa "
$ws.Range("H114").Value = "CRJ9"
$ws.Range("I114").Value = "(D-ACNF)"
$ws.Range("J114").Value = "2:25 PM"
$ws.Range("L114").Value = "0 hours, 55 minutes"

# Row 115
$ws.Range("A115").Value = 114
$ws.Range("B115").Value = "Monday, Jan 09"
$ws.Range("C115").Value = "1:49 PM"
$ws.Range("D115").Value = "P81956"
$ws.Range("E115").Value = "Berlin"
$ws.Range("F115").Value = "(BER)"
$ws.Range("G115").Value = "SprintAir "
$ws.Range("H115").Value = "SF34"
$ws.Range("I115").Value = "(SP-KPV)"
$ws.Range("J115").Value = "1:49 PM"
$ws.Range("L115").Value = "0 hours, 0 minutes"

# Row 116
$ws.Range("A116").Value = 115
$ws.Range("B116").Value = "Monday, Jan 09"
$ws.Range("C116").Value = "2:10 PM"
$ws.Range("D116").Value = "LO3837"
$ws.Range("E116").Value = "Warsaw"
$ws.Range("F116").Value = "(WAW)"
$ws.Range("G116").Value = "LOT "
$ws.Range("H116").Value = "E170"
$ws.Range("I116").Value = "(SP-LDI)"
$ws.Range("J116").Value = "2:08 PM"
$ws.Range("L116").Value = "0 hours, -2 minutes"

# Row 117
$ws.Range("A117").Value = 116
$ws.Range("B117").Value = "Monday, Jan 09"
$ws.Range("C117").Value = "2:40 PM"
$ws.Range("D117").Value = "SK759"
$ws.Range("E117").Value = "Copenhagen"
$ws.Range("F117").Value = "(CPH)"
$ws.Range("G117").Value = "SAS "
$ws.Range("H117").Value = "A20N"
$ws.Range("I117").Value = "(EI-SIA)"
$ws.Range("J117").Value = "2:32 PM"
$ws.Range("L117").Value = "0 hours, -8 minutes"

# Row 118
$ws.Range("A118").Value = 117
$ws.Range("B118").Value = "Monday, Jan 09"
$ws.Range("C118").Value = "2:45 PM"
$ws.Range("D118").Value = "FR6118"
$ws.Range("E118").Value = "London"
$ws.Range("F118").Value = "(STN)"
$ws.Range("G118").Value = "Ryanair "
$ws.Range("H118").Value = "B738"
$ws.Range("I118").Value = "(SP-RSW)"
$ws.Range("J118").Value = "2:39 PM"
$ws.Range("L118").Value = "0 hours, -6 minutes"

# Row 119
$ws.Range("A119").Value = 118
$ws.Range("B119").Value = "Monday, Jan 09"
$ws.Range("C119").Value = "3:20 PM"
$ws.Range("D119").Value = "FR6124"
$ws.Range("E119").Value = "Edinburgh"
$ws.Range("F119").Value = "(EDI)"
$ws.Range("G119").Value = "Ryanair "
$ws.Range("H119").Value = "B738"
$ws.Range("I119").Value = "(SP-RKQ)"
$ws.Range("J119").Value = "3:10 PM"
$ws.Range("L119").Value = "0 hours, -10 minutes"

# Row 120
$ws.Range("A120").Value = 119
$ws.Range("B120").Value = "Monday, Jan 09"
$ws.Range("C120").Value = "3:50 PM"
$ws.Range("D120").Value = "W61732"
$ws.Range("E120").Value = "Stockholm"
$ws.Range("F120").Value = "(NYO)"
$ws.Range("G120").Value = "Wizz Air "
$ws.Range("H120").Value = "A321"
$ws.Range("I120").Value = "(HA-LTB)"
$ws.Range("J120").Value = "3:32 PM"
$ws.Range("L120").Value = "0 hours, -18 minutes"

# Row 121
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = "Monday, Jan 09"
$ws.Range("C121").Value = "4:25 PM"
$ws.Range("D121").Value = "LO3815"
$ws.Range("E121").Value = "Warsaw"
$ws.Range("F121").Value = "(WAW)"
$ws.Range("G121").Value = "LOT "
$ws.Range("H121").Value = "E75S"
$ws.Range("I121").Value = "(SP-LIB)"
$ws.Range("J121").Value = "4:15 PM"
$ws.Range("L121").Value = "0 hours, -10 minutes"

# Row 122
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "Monday, Jan 09"
$ws.Range("C122").Value = "4:42 PM"
$ws.Range("D122").Value = "UNKNOWN"
$ws.Range("E122").Value = "Palanga"
$ws.Range("F122").Value = "(PLQ)"
$ws.Range("G122").Value = "Ryanair "
$ws.Range("H122").Value = "B738"
$ws.Range("I122").Value = "(SP-RSL)"
$ws.Range("J122").Value = "4:53 PM"
$ws.Range("L122").Value = "0 hours, 11 minutes"

# Row 123
$ws.Range("A123").Value = 122
$ws.Range("B123").Value = "Monday, Jan 09"
$ws.Range("C123").Value = "5:30 PM"
$ws.Range("D123").Value = "KL1921"
$ws.Range("E123").Value = "Amsterdam"
$ws.Range("F123").Value = "(AMS)"
$ws.Range("G123").Value = "KLM "
$ws.Range("H123").Value = "E190"
$ws.Range("I123").Value = "(PH-EZG)"
$ws.Range("J123").Value = "5:16 PM"
$ws.Range("L123").Value = "0 hours, -14 minutes"

# Row 124
$ws.Range("A124").Value = 123
$ws.Range("B124").Value = "Monday, Jan 09"
$ws.Range("C124").Value = "5:40 PM"
$ws.Range("D124").Value = "W61784"
$ws.Range("E124").Value = "Oslo"
$ws.Range("F124").Value = "(OSL)"
$ws.Range("G124").Value = "Wizz Air "
$ws.Range("H124").Value = "A320"
$ws.Range("I124").Value = "(HA-LWV)"
$ws.Range("J124").Value = "5:22 PM"
$ws.Range("L124").Value = "0 hours, -18 minutes"

# Row 125
$ws.Range("A125").Value = 124
$ws.Range("B125").Value = "Monday, Jan 09"
$ws.Range("C125").Value = "5:50 PM"
$ws.Range("D125").Value = "W61632"
$ws.Range("E125").Value = "Paris"
$ws.Range("F125").Value = "(BVA)"
$ws.Range("G125").Value = "Wizz Air "
$ws.Range("H125").Value = "A320"
$ws.Range("I125").Value = "(HA-LYS)"
$ws.Range("J125").Value = "5:30 PM"
$ws.Range("L125").Value = "0 hours, -20 minutes"

# Row 126
$ws.Range("A126").Value = 125
$ws.Range("B126").Value = "Monday, Jan 09"
$ws.Range("C126").Value = "5:55 PM"
$ws.Range("D126").Value = "W61746"
$ws.Range("E126").Value = "Bergen"
$ws.Range("F126").Value = "(BGO)"
$ws.Range("G126").Value = "Wizz Air "
$ws.Range("H126").Value = "A321"
$ws.Range("I126").Value = "(HA-LXL)"
$ws.Range("J126").Value = "5:37 PM"
$ws.Range("L126").Value = "0 hours, -18 minutes"
